# Update regression model parameter values across several worksheets.
$wb = $excel.ActiveWorkbook

# regression_model sheet
$ws = $wb.Worksheets.Item("regression_model")
$ws.Range("B2").Value = 0.2360754472511452
$ws.Range("B3").Value = 0.06526393968409917
$ws.Range("B4").Value = 0.1864554834859687
$ws.Range("B5").Value = 0.5016168286296921

# skin_curve sheet
$ws = $wb.Worksheets.Item("skin_curve")
$ws.Range("B2").Value = -2.929461255352919
$ws.Range("B3").Value = 0.07729346639100185
$ws.Range("B4").Value = 0.009714626772967951

# filter_reduction sheet
$ws = $wb.Worksheets.Item("filter_reduction")
$ws.Range("B2").Value = 1.008068653043437
$ws.Range("B3").Value = 0.5342636214790784
$ws.Range("B4").Value = 0.003706813271053893

# fracture_length sheet
$ws = $wb.Worksheets.Item("fracture_length")
$ws.Range("B2").Value = 3.846026467336558
$ws.Range("B3").Value = 0.3181181968665893
$ws.Range("B4").Value = 0.637621805276343
$ws.Range("B5").Value = 0.02693643592064818
